# Generate Report for Handback
# Updates generated timestamps and the "Priority" value for the handback
# status report across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-18 22:15:49"
$overview.Range("G4").Value = "2016-08-18 22:15:49"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E4").Value = "mt"
$zhcn.Range("H2").Value = "2016-08-18 22:15:44"
$zhcn.Range("H4").Value = "2016-08-18 22:15:44"
$zhcn.Range("K2").Value = "2016-08-18 22:16:03"
$zhcn.Range("K4").Value = "2016-08-18 22:16:03"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "mt"
$dede.Range("E4").Value = "mt"
$dede.Range("H2").Value = "2016-08-18 22:15:49"
$dede.Range("H4").Value = "2016-08-18 22:15:49"
$dede.Range("K2").Value = "2016-08-18 22:16:14"
$dede.Range("K4").Value = "2016-08-18 22:16:14"
